# Parameter Passing And Data Provider pgm codes added
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. TableData sheet (sheet index 3): insert a header row above the existing
#    data row, add a new bottom row with the "Tiger Nixon..." summary string.
# ---------------------------------------------------------------------------
$tableData = $wb.Worksheets.Item(3)

# Insert a new row 1; existing row shifts down to row 2.
$tableData.Rows.Item(1).Insert()

# New header row (row 1)
$tableData.Range("A1").Value = "Name"
$tableData.Range("B1").Value = "Position"
$tableData.Range("C1").Value = "Office"
$tableData.Range("D1").Value = "Age"
$tableData.Range("E1").Value = "Start date"
$tableData.Range("F1").Value = "Salary"

# Style the header row: bold, size 12, Times New Roman, color #212529,
# and a medium #DEE2E6 border around every cell.
$headerRng = $tableData.Range("A1:F1")
$headerRng.Font.Name = "Times New Roman"
$headerRng.Font.Size = 12
$headerRng.Font.Bold = $true
$headerRng.Font.Color = 2696481
$headerRng.Borders.Weight = -4138
$headerRng.Borders.Color = 15131358
$tableData.Rows.Item(1).RowHeight = 18.75

# Re-style the (shifted) data row (row 2): same font family/color, but not
# bold, keeping the existing thin border.
$dataRng = $tableData.Range("A2:F2")
$dataRng.Font.Name = "Times New Roman"
$dataRng.Font.Size = 11
$dataRng.Font.Bold = $false
$dataRng.Font.Color = 2696481

# New summary row (row 3)
$tableData.Range("A3").Value = "Tiger Nixon System Architect Edinburgh 61 2011/04/25 `$320,800"

# Column widths
$tableData.Columns.Item(1).ColumnWidth = 57.166666666666664
$tableData.Columns.Item(4).ColumnWidth = 4.022135416666667
$tableData.Columns.Item(6).ColumnWidth = 8.736979166666666

# Selection moves onto the data row now that the header occupies row 1.
$tableData.Range("A2").Select()

# ---------------------------------------------------------------------------
# 2. RadioButtonDemoTestData sheet (sheet index 2): reset the lingering
#    selection back to the top-left cell.
# ---------------------------------------------------------------------------
$radioSheet = $wb.Worksheets.Item(2)
$radioSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 3. New sheet "FormSubmitByParameter" holding the submit-confirmation data.
# ---------------------------------------------------------------------------
$formSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$formSheet.Name = "FormSubmitByParameter"

$formSheet.Range("A1").Value = "Expected Message"
$formSheet.Range("A1").Font.Bold = $true
$formSheet.Range("A1").HorizontalAlignment = -4131
$formSheet.Range("A1").VerticalAlignment = -4160

$formSheet.Range("A2").Value = "Form has been submitted successfully!"
$formSheet.Range("A2").Font.Name = "Times New Roman"
$formSheet.Range("A2").Font.Size = 11
$formSheet.Range("A2").Font.Bold = $false
$formSheet.Range("A2").Font.Color = 2696481

$formSheet.Columns.Item(1).ColumnWidth = 43.44401041666667

$formSheet.Range("A1:A2").Select()
